$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O2").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("O3").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("O4").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("O5").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("O6").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("O7").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("O8").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("O9").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("O10").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("O11").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("O12").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("O13").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("O14").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("O15").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("O16").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("O17").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("O18").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("O19").Value = "de503c24-f17d-47a9-9a47-6f0a194f8c9c"
$ws.Range("O20").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("O21").Value = "718c6b8f-7c00-4bcb-b53c-8f3f42154362"
$ws.Range("O22").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("O23").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("O24").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
